$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (26) since dimension shrinks from 26 to 25 rows
$ws.Rows(26).Delete()

$arr = New-Object 'object[,]' 24,9
$arr[0,0] = "model_9_3_12"
$arr[0,1] = -0.4815702885590147
$arr[0,2] = -15.99720160848214
$arr[0,3] = -10.04781711999013
$arr[0,4] = -11.83924151271215
$arr[0,5] = 1.63966178894043
$arr[0,6] = 5.129796981811523
$arr[0,7] = 5.164009094238281
$arr[0,8] = 5.145896911621094

$arr[1,0] = "model_9_3_11"
$arr[1,1] = -0.4596741393442021
$arr[1,2] = -17.02915414119099
$arr[1,3] = -8.650423278569177
$arr[1,4] = -11.48371315024235
$arr[1,5] = 1.615429162979126
$arr[1,6] = 5.441243171691895
$arr[1,7] = 4.510834693908691
$arr[1,8] = 5.003402709960938

$arr[2,0] = "model_9_3_10"
$arr[2,1] = -0.4470424457085991
$arr[2,2] = -17.60570524260995
$arr[2,3] = -7.892471760701275
$arr[2,4] = -11.29757912506017
$arr[2,5] = 1.601449608802795
$arr[2,6] = 5.61524772644043
$arr[2,7] = 4.15654993057251
$arr[2,8] = 4.928801536560059

$arr[3,0] = "model_9_3_9"
$arr[3,1] = -0.1987306924556365
$arr[3,2] = -13.55199367819378
$arr[3,3] = -6.657028950360103
$arr[3,4] = -9.003521463305455
$arr[3,5] = 1.32664167881012
$arr[3,6] = 4.391827583312988
$arr[3,7] = 3.57907509803772
$arr[3,8] = 4.009356021881104

$arr[4,0] = "model_9_3_8"
$arr[4,1] = -0.1289760336038486
$arr[4,2] = -12.31325148276143
$arr[4,3] = -6.285733118579742
$arr[4,4] = -8.305918316872335
$arr[4,5] = 1.249443888664246
$arr[4,6] = 4.017971515655518
$arr[4,7] = 3.40552282333374
$arr[4,8] = 3.729760885238647

$arr[5,0] = "model_9_3_1"
$arr[5,1] = -0.1075053005549194
$arr[5,2] = -9.65189943724573
$arr[5,3] = -4.078260566676644
$arr[5,4] = -6.033458361419521
$arr[5,5] = 1.225682020187378
$arr[5,6] = 3.214769601821899
$arr[5,7] = 2.373698234558105
$arr[5,8] = 2.818971395492554

$arr[6,0] = "model_9_3_22"
$arr[6,1] = -0.08531816016939442
$arr[6,2] = -12.21307573077441
$arr[6,3] = -7.25751964104613
$arr[6,4] = -8.799321570937655
$arr[6,5] = 1.201127290725708
$arr[6,6] = 3.987738609313965
$arr[6,7] = 3.859758377075195
$arr[6,8] = 3.92751407623291

$arr[7,0] = "model_9_3_21"
$arr[7,1] = -0.08143757889779324
$arr[7,2] = -12.22165177946186
$arr[7,3] = -7.15173724520813
$arr[7,4] = -8.744683513920869
$arr[7,5] = 1.196832656860352
$arr[7,6] = 3.990326881408691
$arr[7,7] = 3.81031322479248
$arr[7,8] = 3.905615329742432

$arr[8,0] = "model_9_3_23"
$arr[8,1] = -0.08020733176338002
$arr[8,2] = -12.18883391163021
$arr[8,3] = -7.213352206133736
$arr[8,4] = -8.76541414280036
$arr[8,5] = 1.195471167564392
$arr[8,6] = 3.980422496795654
$arr[8,7] = 3.839113712310791
$arr[8,8] = 3.913923740386963

$arr[9,0] = "model_9_3_18"
$arr[9,1] = -0.07750576280355759
$arr[9,2] = -11.24811659983113
$arr[9,3] = -7.658826208379734
$arr[9,4] = -8.634881628396727
$arr[9,5] = 1.192481279373169
$arr[9,6] = 3.696511507034302
$arr[9,7] = 4.047338485717773
$arr[9,8] = 3.861607313156128

$arr[10,0] = "model_9_3_20"
$arr[10,1] = -0.07246103768415368
$arr[10,2] = -11.56533692963427
$arr[10,3] = -7.47091012383104
$arr[10,4] = -8.658213722510562
$arr[10,5] = 1.186898350715637
$arr[10,6] = 3.792249441146851
$arr[10,7] = 3.959502220153809
$arr[10,8] = 3.870958566665649

$arr[11,0] = "model_9_3_13"
$arr[11,1] = -0.06527118762241058
$arr[11,2] = -10.34297128291709
$arr[11,3] = -7.79982437057752
$arr[11,4] = -8.351422833027877
$arr[11,5] = 1.178941369056702
$arr[11,6] = 3.423336744308472
$arr[11,7] = 4.113244533538818
$arr[11,8] = 3.747998476028442

$arr[12,0] = "model_9_3_14"
$arr[12,1] = -0.05672633646216574
$arr[12,2] = -10.16030449279204
$arr[12,3] = -7.823599885829825
$arr[12,4] = -8.291650437899188
$arr[12,5] = 1.169484734535217
$arr[12,6] = 3.368206977844238
$arr[12,7] = 4.1243577003479
$arr[12,8] = 3.724042177200317

$arr[13,0] = "model_9_3_15"
$arr[13,1] = -0.05513473165281835
$arr[13,2] = -10.32213043405362
$arr[13,3] = -7.761754724044282
$arr[13,4] = -8.322221374085569
$arr[13,5] = 1.167723298072815
$arr[13,6] = 3.417046546936035
$arr[13,7] = 4.095449924468994
$arr[13,8] = 3.736294984817505

$arr[14,0] = "model_9_3_17"
$arr[14,1] = -0.05507444295549524
$arr[14,2] = -10.36031082621072
$arr[14,3] = -7.79830500027084
$arr[14,4] = -8.357500015674322
$arr[14,5] = 1.167656540870667
$arr[14,6] = 3.428569793701172
$arr[14,7] = 4.112534523010254
$arr[14,8] = 3.750433921813965

$arr[15,0] = "model_9_3_16"
$arr[15,1] = -0.05503169369445415
$arr[15,2] = -10.38705592081347
$arr[15,3] = -7.756997988770578
$arr[15,4] = -8.345490658553175
$arr[15,5] = 1.167609214782715
$arr[15,6] = 3.436641216278076
$arr[15,7] = 4.093226432800293
$arr[15,8] = 3.745621204376221

$arr[16,0] = "model_9_3_19"
$arr[16,1] = -0.0530942086714723
$arr[16,2] = -10.94162774139916
$arr[16,3] = -7.53878941091344
$arr[16,4] = -8.446817885938296
$arr[16,5] = 1.165464997291565
$arr[16,6] = 3.604012489318848
$arr[16,7] = 3.991230487823486
$arr[16,8] = 3.786232471466064

$arr[17,0] = "model_9_3_7"
$arr[17,1] = -0.04108163239206974
$arr[17,2] = -10.80183556860581
$arr[17,3] = -5.691914686328423
$arr[17,4] = -7.377489101036756
$arr[17,5] = 1.152170538902283
$arr[17,6] = 3.561823129653931
$arr[17,7] = 3.127958059310913
$arr[17,8] = 3.357651233673096

$arr[18,0] = "model_9_3_6"
$arr[18,1] = -0.02257449890778496
$arr[18,2] = -10.51001009701293
$arr[18,3] = -5.361365681931057
$arr[18,4] = -7.079740132911024
$arr[18,5] = 1.131688714027405
$arr[18,6] = 3.473749160766602
$arr[18,7] = 2.973451614379883
$arr[18,8] = 3.238315343856812

$arr[19,0] = "model_9_3_5"
$arr[19,1] = 0.05106327205737293
$arr[19,2] = -10.14789594023329
$arr[19,3] = -4.056410488568249
$arr[19,4] = -6.219196892755166
$arr[19,5] = 1.050193428993225
$arr[19,6] = 3.364462375640869
$arr[19,7] = 2.363484859466553
$arr[19,8] = 2.893414258956909

$arr[20,0] = "model_9_3_4"
$arr[20,1] = 0.07919627605960333
$arr[20,2] = -9.917223758639839
$arr[20,3] = -3.676986065496357
$arr[20,4] = -5.919003041120866
$arr[20,5] = 1.019058346748352
$arr[20,6] = 3.294844627380371
$arr[20,7] = 2.186133146286011
$arr[20,8] = 2.773098230361938

$arr[21,0] = "model_9_3_2"
$arr[21,1] = 0.2223464636237602
$arr[21,2] = -7.679985552177822
$arr[21,3] = -2.234330297006648
$arr[21,4] = -4.235363601288164
$arr[21,5] = 0.8606333136558533
$arr[21,6] = 2.619641065597534
$arr[21,7] = 1.511801838874817
$arr[21,8] = 2.098304748535156

$arr[22,0] = "model_9_3_3"
$arr[22,1] = 0.2589132038618094
$arr[22,2] = -7.226364267672183
$arr[22,3] = -2.416041297212084
$arr[22,4] = -4.154253816728081
$arr[22,5] = 0.820164680480957
$arr[22,6] = 2.482736825942993
$arr[22,7] = 1.596737861633301
$arr[22,8] = 2.065796375274658

$arr[23,0] = "model_9_3_0"
$arr[23,1] = 0.3583886394316547
$arr[23,2] = -1.579100324131243
$arr[23,3] = -1.148739557365386
$arr[23,4] = -1.207436167593781
$arr[23,5] = 0.7100747227668762
$arr[23,6] = 0.7783787846565247
$arr[23,7] = 1.004371285438538
$arr[23,8] = 0.8847281336784363

$ws.Range("A2:I25").Value = $arr